$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.769.58'
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").Value = '3.109.40'

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.06'
$ws.Range("E5").Value = '  -3.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '613.53'
$ws.Range("E6").Value = '  -0.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.08'
$ws.Range("E7").Value = '  -1.73%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.389'
$ws.Range("E8").Value = '  -0.65%  '

$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").Value = '3.107.49'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.786'
$ws.Range("E11").Value = '  +3.90%  '

$ws.Range("E12").Value = '  -3.61%  '

$ws.Range("E13").Value = '  -3.45%  '

$ws.Range("D14").Value = '92.504.67'
$ws.Range("E14").Value = '  +0.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.98'
$ws.Range("E15").Value = '  -3.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.42'
$ws.Range("E16").Value = '  -3.18%  '

$ws.Range("D18").Value = '3.117.05'
$ws.Range("E18").Value = '  -0.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.83'
$ws.Range("E19").Value = '  +1.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.47'
$ws.Range("E20").Value = '  -3.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.83'
$ws.Range("E21").Value = '  -0.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000203'
$ws.Range("E22").Value = '  -0.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '439.42'
$ws.Range("E23").Value = '  -3.62%  '

$ws.Range("E24").Value = '  -0.80%  '

$ws.Range("E25").Value = '  -5.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '85.44'
$ws.Range("E26").Value = '  -4.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.74'
$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").Value = '3.274.39'
$ws.Range("E28").Value = '  -0.67%  '

$ws.Range("E29").Value = '  +0.37%  '

$ws.Range("E30").Value = '  +8.22%  '

$ws.Range("E31").Value = '  -9.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.233'
$ws.Range("E32").Value = '  +2.41%  '

$ws.Range("E33").Value = '  -30.35%  '

$ws.Range("E34").Value = '  -2.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.12'
$ws.Range("E35").Value = '  +8.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.162'
$ws.Range("E36").Value = '  -6.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.69'
$ws.Range("E37").Value = '  -2.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.99'
$ws.Range("E38").Value = '  +3.96%  '

$ws.Range("E39").Value = '  -8.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.93'
$ws.Range("E40").Value = '  +7.76%  '

$ws.Range("E41").Value = '  -2.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '466.48'
$ws.Range("E42").Value = '  -5.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.427'
$ws.Range("E43").Value = '  -2.33%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.35'
$ws.Range("E44").Value = '  -2.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '159.04'
$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("E47").Value = '  -3.47%  '

$ws.Range("E48").Value = '  -4.97%  '

$ws.Range("E49").Value = '  -2.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.81'
$ws.Range("E50").Value = '  -0.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0323'
$ws.Range("E51").Value = '  -0.58%  '
